$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Data for the three new columns (I: defect ID, J: date, K: priority)
# mirroring the style/format used by columns A (integer) and B (date)
$ids       = @(41, 42, 43, 44, 45, 46)
$dates     = @(42005, 42007, 42005, 42007, 42005, 42007)
$priority  = @(1, 2, 1, 2, 1, 2)

for ($i = 0; $i -lt 6; $i++) {
    $r = 7 + $i

    $cellI = $ws.Cells.Item($r, 9)   # column I
    $cellI.Value = $ids[$i]
    $cellI.NumberFormat = "0"

    $cellJ = $ws.Cells.Item($r, 10)  # column J
    $cellJ.Value = $dates[$i]
    $cellJ.NumberFormat = [char]91 + "`$-409]d\-mmm\-yy;@"

    $cellK = $ws.Cells.Item($r, 11)  # column K
    $cellK.Value = $priority[$i]
}

# Update the selection to match the new active range
$ws.Range("I7:K12").Select()
